# Remove the trailing "Análisis costo-beneficio" slide (last slide in the
# deck). The cost/benefit analysis now lives in the separate PDF document
# that accompanies the proposal, so it is dropped from the sldIdLst here.
$p = $ppt.ActivePresentation
$lastIndex = $p.Slides.Count
$p.Slides.Item($lastIndex).Delete()
